$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so numeric-looking strings (e.g. "1.00") are not
# auto-converted to numbers, then restore the default style so no stray
# cell-style index is introduced.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '62.095.16'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '3.017.75'
$ws.Range("E3").Value = '  -1.22%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").Value = '541.11'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '132.40'
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '3.007.90'
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").Value = '6.10'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = '0.146'
$ws.Range("E11").Value = '  -4.51%  '
$ws.Range("D12").Value = '0.443'
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '33.89'
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").Value = '3.523.95'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").Value = '62.214.99'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '3.029.98'
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("E18").Value = '  -3.78%  '
$ws.Range("D19").Value = '6.58'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '475.19'
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").Value = '13.18'
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").Value = '0.668'
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("D23").Value = '7.00'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '80.72'
$ws.Range("E24").Value = '  +3.19%  '
$ws.Range("D25").Value = '12.02'
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").Value = '7.71'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("D30").Value = '1.92'
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").Value = '25.53'
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").Value = '1.12'
$ws.Range("E32").Value = '  -2.53%  '
$ws.Range("D33").Value = '2.34'
$ws.Range("E33").Value = '  +2.81%  '
$ws.Range("D34").Value = '5.59'
$ws.Range("E34").Value = '  +2.49%  '
$ws.Range("D35").Value = '54.68'
$ws.Range("E35").Value = '  -5.10%  '
$ws.Range("D36").Value = '5.82'
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = '457.19'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").Value = '3.154.21'
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").Value = '0.0796'
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("D40").Value = '0.0385'
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("D41").Value = '0.117'
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("D42").Value = '8.04'
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  -3.73%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '26.18'
$ws.Range("E45").Value = '  +3.84%  '
$ws.Range("D46").Value = '0.242'
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("D47").Value = '0.108'
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = '1.96'
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '114.17'
$ws.Range("E49").Value = '  -6.35%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0496'
$ws.Range("E50").Value = '  -3.78%  '
$ws.Range("E51").Value = '  +3.05%  '

$ws.Range("D2:E51").Style = "Normal"
